$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5333
$ws.Range("C3").Value = 9615
$ws.Range("D3").Value = 17100
$ws.Range("E3").Value = 24500
$ws.Range("F3").Value = 24200
$ws.Range("G3").Value = 21300
$ws.Range("B4").Value = 699.4001919999999
$ws.Range("C4").Value = 1260.388352
$ws.Range("D4").Value = 2246.049792
$ws.Range("E4").Value = 3205.496832
$ws.Range("F4").Value = 3177.18528
$ws.Range("G4").Value = 2789.21216
$ws.Range("B5").Value = 186.87
$ws.Range("C5").Value = 196.6
$ws.Range("D5").Value = 225.29
$ws.Range("E5").Value = 314.86
$ws.Range("F5").Value = 646.47
$ws.Range("G5").Value = 1475.04
$ws.Range("B6").Value = 310
$ws.Range("C6").Value = 351
$ws.Range("D6").Value = 396
$ws.Range("E6").Value = 586
$ws.Range("F6").Value = 1287
$ws.Range("G6").Value = 2999
$ws.Range("B7").Value = 363
$ws.Range("C7").Value = 388
$ws.Range("D7").Value = 461
$ws.Range("F7").Value = 1500
$ws.Range("G7").Value = 3720
$ws.Range("B12").Value = 25700
$ws.Range("C12").Value = 51600
$ws.Range("D12").Value = 100000
$ws.Range("E12").Value = 189000
$ws.Range("F12").Value = 335000
$ws.Range("G12").Value = 503000
$ws.Range("B13").Value = 104.8576
$ws.Range("C13").Value = 210.763776
$ws.Range("D13").Value = 411.041792
$ws.Range("E13").Value = 773.8490880000001
$ws.Range("F13").Value = 1371.537408
$ws.Range("G13").Value = 2060.45184
$ws.Range("B14").Value = 38.61556
$ws.Range("C14").Value = 37.99648000000001
$ws.Range("D14").Value = 39.25501000000001
$ws.Range("E14").Value = 41.65730000000001
$ws.Range("F14").Value = 46.75563
$ws.Range("G14").Value = 62.08248
$ws.Range("B15").Value = 72.19200000000001
$ws.Range("C15").Value = 69.12
$ws.Range("D15").Value = 75.264
$ws.Range("E15").Value = 84.48
$ws.Range("F15").Value = 98.816
$ws.Range("G15").Value = 134.144
$ws.Range("B16").Value = 94.72
$ws.Range("C16").Value = 96.768
$ws.Range("D16").Value = 99.84
$ws.Range("E16").Value = 108.032
$ws.Range("F16").Value = 126.464
$ws.Range("G16").Value = 177.152
$ws.Range("B21").Value = 7111
$ws.Range("C21").Value = 11500
$ws.Range("D21").Value = 17500
$ws.Range("E21").Value = 20100
$ws.Range("F21").Value = 19500
$ws.Range("G21").Value = 15200
$ws.Range("B22").Value = 932.184064
$ws.Range("C22").Value = 1507.852288
$ws.Range("D22").Value = 2294.284288
$ws.Range("E22").Value = 2638.217216
$ws.Range("F22").Value = 2556.428288
$ws.Range("G22").Value = 1987.05152
$ws.Range("B23").Value = 105.81
$ws.Range("C23").Value = 109.74
$ws.Range("D23").Value = 123.61
$ws.Range("E23").Value = 186.93
$ws.Range("F23").Value = 464.54
$ws.Range("G23").Value = 1417.81
$ws.Range("B24").Value = 121
$ws.Range("C24").Value = 135
$ws.Range("D24").Value = 184
$ws.Range("E24").Value = 289
$ws.Range("F24").Value = 963
$ws.Range("G24").Value = 3621
$ws.Range("B25").Value = 167
$ws.Range("C25").Value = 174
$ws.Range("D25").Value = 221
$ws.Range("E25").Value = 355
$ws.Range("F25").Value = 1516
$ws.Range("G25").Value = 4555
$ws.Range("B30").Value = 91300
$ws.Range("C30").Value = 147000
$ws.Range("D30").Value = 278000
$ws.Range("E30").Value = 425000
$ws.Range("F30").Value = 459000
$ws.Range("G30").Value = 379000
$ws.Range("B31").Value = 374.341632
$ws.Range("C31").Value = 601.882624
$ws.Range("D31").Value = 1137.70496
$ws.Range("E31").Value = 1740.63616
$ws.Range("F31").Value = 1882.19392
$ws.Range("G31").Value = 1552.941056
$ws.Range("B32").Value = 9.121370000000001
$ws.Range("C32").Value = 11.16
$ws.Range("D32").Value = 10.49798
$ws.Range("E32").Value = 11.54
$ws.Range("F32").Value = 15.64
$ws.Range("G32").Value = 46.21
$ws.Range("B33").Value = 10.944
$ws.Range("C33").Value = 13.888
$ws.Range("D33").Value = 13.376
$ws.Range("E33").Value = 16.32
$ws.Range("F33").Value = 21.888
$ws.Range("G33").Value = 69
$ws.Range("B34").Value = 12.608
$ws.Range("C34").Value = 16.32
$ws.Range("D34").Value = 15.68
$ws.Range("E34").Value = 20.352
$ws.Range("F34").Value = 29.568
$ws.Range("G34").Value = 88
$ws.Range("B39").Value = 10300
$ws.Range("C39").Value = 9846
$ws.Range("D39").Value = 15900
$ws.Range("E39").Value = 13000
$ws.Range("G39").Value = 10900
$ws.Range("B40").Value = 1355.808768
$ws.Range("C40").Value = 1290.797056
$ws.Range("D40").Value = 2088.763392
$ws.Range("E40").Value = 1831.862272
$ws.Range("F40").Value = 1665.138688
$ws.Range("G40").Value = 1430.257664
$ws.Range("B41").Value = 95.94
$ws.Range("C41").Value = 160.39
$ws.Range("D41").Value = 230.78
$ws.Range("E41").Value = 556.4299999999999
$ws.Range("F41").Value = 1231.82
$ws.Range("G41").Value = 2858.32
$ws.Range("B42").Value = 225
$ws.Range("C42").Value = 310
$ws.Range("D42").Value = 619
$ws.Range("E42").Value = 1434
$ws.Range("F42").Value = 3621
$ws.Range("G42").Value = 9765
$ws.Range("B43").Value = 273
$ws.Range("C43").Value = 502
$ws.Range("D43").Value = 840
$ws.Range("E43").Value = 1893
$ws.Range("F43").Value = 5342
$ws.Range("G43").Value = 18220
$ws.Range("B48").Value = 309000
$ws.Range("C48").Value = 455000
$ws.Range("D48").Value = 520000
$ws.Range("E48").Value = 461000
$ws.Range("F48").Value = 429000
$ws.Range("G48").Value = 383000
$ws.Range("B49").Value = 1266.679808
$ws.Range("C49").Value = 1864.368128
$ws.Range("D49").Value = 2130.706432
$ws.Range("E49").Value = 1887.4368
$ws.Range("F49").Value = 1757.413376
$ws.Range("G49").Value = 1566.572544
$ws.Range("B50").Value = 2.97914
$ws.Range("C50").Value = 3.90883
$ws.Range("D50").Value = 7.23163
$ws.Range("E50").Value = 16.57516
$ws.Range("F50").Value = 35.39712
$ws.Range("G50").Value = 79.96814000000001
$ws.Range("B51").Value = 1.608
$ws.Range("C51").Value = 2.096
$ws.Range("D51").Value = 4.96
$ws.Range("E51").Value = 10.688
$ws.Range("F51").Value = 24.96
$ws.Range("G51").Value = 31.616
$ws.Range("B52").Value = 102.912
$ws.Range("C52").Value = 138.24
$ws.Range("D52").Value = 232.448
$ws.Range("E52").Value = 577.5360000000001
$ws.Range("F52").Value = 1122.304
$ws.Range("G52").Value = 2179.072
$ws.Range("B57").Value = 5953
$ws.Range("C57").Value = 8827
$ws.Range("D57").Value = 13500
$ws.Range("E57").Value = 15300
$ws.Range("F57").Value = 16100
$ws.Range("G57").Value = 14500
$ws.Range("B58").Value = 780.140544
$ws.Range("C58").Value = 1156.579328
$ws.Range("D58").Value = 1765.801984
$ws.Range("E58").Value = 2011.168768
$ws.Range("F58").Value = 2109.734912
$ws.Range("G58").Value = 1897.92256
$ws.Range("B59").Value = 116.6
$ws.Range("C59").Value = 119.87
$ws.Range("D59").Value = 125.42
$ws.Range("E59").Value = 176.36
$ws.Range("F59").Value = 494.11
$ws.Range("G59").Value = 1497.74
$ws.Range("B60").Value = 124
$ws.Range("C60").Value = 135
$ws.Range("D60").Value = 147
$ws.Range("E60").Value = 247
$ws.Range("F60").Value = 865
$ws.Range("G60").Value = 3392
$ws.Range("B61").Value = 149
$ws.Range("C61").Value = 161
$ws.Range("D61").Value = 182
$ws.Range("E61").Value = 322
$ws.Range("F61").Value = 1188
$ws.Range("G61").Value = 3851
$ws.Range("B66").Value = 91300
$ws.Range("C66").Value = 163000
$ws.Range("D66").Value = 254000
$ws.Range("E66").Value = 370000
$ws.Range("F66").Value = 497000
$ws.Range("G66").Value = 474000
$ws.Range("B67").Value = 374.341632
$ws.Range("C67").Value = 668.991488
$ws.Range("D67").Value = 1040.187392
$ws.Range("E67").Value = 1514.143744
$ws.Range("F67").Value = 2035.286016
$ws.Range("G67").Value = 1940.914176
$ws.Range("B68").Value = 9.06779
$ws.Range("C68").Value = $null
$ws.Range("D68").Value = 10.67
$ws.Range("E68").Value = 11.15416
$ws.Range("F68").Value = 16.42
$ws.Range("G68").Value = 50.36
$ws.Range("B69").Value = 10.176
$ws.Range("C69").Value = 10.816
$ws.Range("D69").Value = 14.272
$ws.Range("E69").Value = 15.68
$ws.Range("F69").Value = 24
$ws.Range("G69").Value = 81
$ws.Range("B70").Value = 11.328
$ws.Range("C70").Value = 12.352
$ws.Range("D70").Value = 16.768
$ws.Range("E70").Value = 19.328
$ws.Range("F70").Value = 35
$ws.Range("G70").Value = 110
